$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("create")

# --- Header row (row 1): new "loyalties" triple (name/expected/message style) ---
$ws.Range("J1").Value = "loyalties"
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("K1").Value = "expected"
$ws.Range("B1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

$ws.Range("L1").Value = "message"
$ws.Range("C1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

# --- Data row (row 2) ---
$ws.Range("J2").Value = '[{"id":1,"pointBalance":2},{"id":2,"pointBalance":3}]'
$ws.Range("A3").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").WrapText = $true

$ws.Range("K2").Value = $true
$ws.Range("A3").Copy()
$ws.Range("K2").PasteSpecial(-4122)

$ws.Range("L2").Value = "Thêm mới thành công."
$ws.Range("A3").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column widths ---
$ws.Columns.Item(10).ColumnWidth = 19.6640625
$ws.Columns.Item(12).ColumnWidth = 19.33203125
$ws.Columns.Item(29).ColumnWidth = 17.5546875

# --- View state ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I10").Select()

Write-Output "done"
